$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage (preserves leading zeros /
# numeric-looking IDs like case numbers, OT numbers and dates-as-text), then
# strip the temporary "@" number-format style so the cell is left exactly as
# it was originally (no explicit style index).
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value() = $val
    $rng.ClearFormats()
}

# The automated map refresh dropped one resolved case from the bottom of the
# table (old row 89, case 7060/PINZON 1578 becomes row 87's new content) and
# the remaining rows shift their content up by one; row 88 receives the
# brand-new case reported on 8/26/2025 (Ameghino 523). Implement this as: a
# structural delete of the last row, followed by writing each row's new
# content directly (matches how the source data feed rewrote the sheet).
$ws.Rows(89).Delete()

Set-TextValue $ws.Range("A80") "6561"
Set-TextValue $ws.Range("C80") "GOLETA SARANDI 6050 "
Set-TextValue $ws.Range("D80") "8"
Set-TextValue $ws.Range("E80") "808733912"
Set-TextValue $ws.Range("H80") "Picada"
Set-TextValue $ws.Range("J80") "Cambio"
$ws.Range("M80").Value() = -58.468841
$ws.Range("N80").Value() = -34.686635
Set-TextValue $ws.Range("A81") "6567"
Set-TextValue $ws.Range("C81") "SAN ANTONIO 1221"
Set-TextValue $ws.Range("D81") "4"
Set-TextValue $ws.Range("E81") "808733914"
Set-TextValue $ws.Range("H81") "Aplomar"
Set-TextValue $ws.Range("J81") "Aplomo"
$ws.Range("M81").Value() = -58.375684
$ws.Range("N81").Value() = -34.656092
Set-TextValue $ws.Range("O81") "San Telmo"
Set-TextValue $ws.Range("A82") "-549"
Set-TextValue $ws.Range("B82") "8/7/2025"
Set-TextValue $ws.Range("C82") "14 de Julio 65"
Set-TextValue $ws.Range("D82") "13"
Set-TextValue $ws.Range("E82") "808749189"
Set-TextValue $ws.Range("H82") "picada"
Set-TextValue $ws.Range("J82") "Cambio"
$ws.Range("M82").Value() = -58.468496
$ws.Range("N82").Value() = -34.591282
Set-TextValue $ws.Range("O82") "Paternal"
Set-TextValue $ws.Range("P82") "Capital Norte"
Set-TextValue $ws.Range("A83") "6960"
Set-TextValue $ws.Range("B83") "8/14/2025"
Set-TextValue $ws.Range("C83") "VALLESE, FELIPE 1940"
Set-TextValue $ws.Range("D83") "7"
Set-TextValue $ws.Range("E83") "808972988"
Set-TextValue $ws.Range("H83") "Picada"
$ws.Range("M83").Value() = -58.460818
$ws.Range("N83").Value() = -34.618934
Set-TextValue $ws.Range("O83") "Boedo"
Set-TextValue $ws.Range("P83") "Capital Sur"
Set-TextValue $ws.Range("A84") "6979"
Set-TextValue $ws.Range("B84") "8/18/2025"
Set-TextValue $ws.Range("C84") "RIVADAVIA AV. 6740"
Set-TextValue $ws.Range("D84") "7"
Set-TextValue $ws.Range("E84") "809006419"
Set-TextValue $ws.Range("H84") "Reclaman fuera de plomo ver si es necesario cambio"
Set-TextValue $ws.Range("J84") "Aplomo"
$ws.Range("M84").Value() = -58.460441
$ws.Range("N84").Value() = -34.628243
Set-TextValue $ws.Range("O84") "Boedo"
Set-TextValue $ws.Range("A85") "-557"
Set-TextValue $ws.Range("B85") "8/21/2025"
Set-TextValue $ws.Range("C85") "Av Castañares 4621"
Set-TextValue $ws.Range("D85") "8"
Set-TextValue $ws.Range("E85") "ICD30462144"
Set-TextValue $ws.Range("H85") "Colocar columna para pedir traspaso de nodo telecom"
Set-TextValue $ws.Range("K85") "Nodo Teco"
$ws.Range("M85").Value() = -58.470977
$ws.Range("N85").Value() = -34.665358
Set-TextValue $ws.Range("A86") "7051"
Set-TextValue $ws.Range("B86") "8/26/2025"
Set-TextValue $ws.Range("C86") "MORENO, JOSE MARIA AV. 345"
Set-TextValue $ws.Range("D86") "6"
Set-TextValue $ws.Range("E86") "ICD30508311"
Set-TextValue $ws.Range("H86") "Colocar PRFV R400 para pedir traspaso de fuente"
Set-TextValue $ws.Range("J86") "Cambio"
Set-TextValue $ws.Range("K86") "Fuente Teco"
Set-TextValue $ws.Range("L86") "Terminal"
$ws.Range("M86").Value() = -58.435017
$ws.Range("N86").Value() = -34.622044
Set-TextValue $ws.Range("A87") "7060"
Set-TextValue $ws.Range("B87") "8/26/2025"
Set-TextValue $ws.Range("C87") "PINZON 1578"
Set-TextValue $ws.Range("D87") "4"
Set-TextValue $ws.Range("E87") "809195671"
Set-TextValue $ws.Range("H87") "Picada"
Set-TextValue $ws.Range("K87") "Sin equipos"
$ws.Range("M87").Value() = -58.373428
$ws.Range("N87").Value() = -34.63705
Set-TextValue $ws.Range("O87") "San Telmo"
Set-TextValue $ws.Range("A88") "-568"
Set-TextValue $ws.Range("C88") "Ameghino 523"
Set-TextValue $ws.Range("D88") "9"
Set-TextValue $ws.Range("E88") "809208239"
Set-TextValue $ws.Range("H88") "Picada con pelgro de caida"
Set-TextValue $ws.Range("K88") "Sin equipos"
Set-TextValue $ws.Range("L88") "Pasante"
$ws.Range("M88").Value() = -58.488424
$ws.Range("N88").Value() = -34.642002
Set-TextValue $ws.Range("O88") "Devoto"
Set-TextValue $ws.Range("P88") "Capital Norte"

Write-Host "Done updating AYKO sheet: dimension now" $ws.UsedRange.Address()
